# "Update UI Appframework 3"
#
# The "Encrypt" sheet (3rd tab) used to hold two plain-text DES/MD5 notes
# in B3/B5. They're replaced with two hyperlinked reference URLs (font-face
# / icon-font docs) written into A2/A4, and that sheet becomes the active
# tab. The "Server Auth Rule" sheet (2nd tab) just loses its "active" flag
# and scrolls down a bit (no value changes other than the shared-string
# reindex caused by removing the two old strings, which Excel does for us
# automatically).

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)   # "Server Auth Rule"
$ws3 = $wb.Worksheets.Item(3)   # "Encrypt"

# --- Encrypt sheet: drop the old DES/MD5 text cells ---------------------
$ws3.Range("B3").ClearContents()
$ws3.Range("B5").ClearContents()

# --- Encrypt sheet: add the two new hyperlinked reference cells ---------
$ws3.Hyperlinks.Add($ws3.Range("A2"), "http://www.w3cplus.com/content/css3-font-face")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://icomoon.io/app/#/select", "/select")

# --- Server Auth Rule sheet: scroll so row 13 is near the top -----------
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 13

# --- Encrypt sheet becomes the active tab, with K10 selected ------------
$ws3.Activate()
$ws3.Range("K10").Select() | Out-Null
